$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D:G columns (rows 2-51) to keep their existing "Text" cell type instead of
# Excel auto-inferring numeric-looking strings ("307.35", "1.54%", "3", ...) as numbers.
$numericLookingRange = $ws.Range("D2:G51")
$numericLookingRange.NumberFormat = "@"

$ws.Range("D2").Value = '307.35'
$ws.Range("E2").Value = '1.54%'
$ws.Range("G2").Value = '3'
$ws.Range("D3").Value = '39.23'
$ws.Range("E3").Value = '10.32%'
$ws.Range("G3").Value = '3'
$ws.Range("E4").Value = '0.89%'
$ws.Range("G4").Value = '3'
$ws.Range("D5").Value = '0.08146'
$ws.Range("E5").Value = '3.09%'
$ws.Range("G5").Value = '3'
$ws.Range("D6").Value = '1.970'
$ws.Range("E6").Value = '6.73%'
$ws.Range("G6").Value = '3'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '4.176'
$ws.Range("E7").Value = '1.79%'
$ws.Range("G7").Value = '3'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").Value = '7.899'
$ws.Range("E8").Value = '1.52%'
$ws.Range("G8").Value = '3'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9288'
$ws.Range("E9").Value = '0.94%'
$ws.Range("G9").Value = '3'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1419'
$ws.Range("E10").Value = '5.26%'
$ws.Range("G10").Value = '3'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1946'
$ws.Range("E11").Value = '2.32%'
$ws.Range("G11").Value = '3'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.09250'
$ws.Range("E12").Value = '1.31%'
$ws.Range("G12").Value = '3'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03504'
$ws.Range("E13").Value = '1.19%'
$ws.Range("G13").Value = '3'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09876'
$ws.Range("E14").Value = '0.42%'
$ws.Range("G14").Value = '3'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001421'
$ws.Range("E15").Value = '0.91%'
$ws.Range("G15").Value = '3'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005999'
$ws.Range("E16").Value = '-1.39%'
$ws.Range("G16").Value = '3'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.946'
$ws.Range("E17").Value = '5.85%'
$ws.Range("G17").Value = '3'
$ws.Range("D18").Value = '3.454'
$ws.Range("E18").Value = '3.43%'
$ws.Range("G18").Value = '3'
$ws.Range("D19").Value = '0.3453'
$ws.Range("G19").Value = '3'
$ws.Range("D20").Value = '0.1328'
$ws.Range("E20").Value = '-1.15%'
$ws.Range("G20").Value = '3'
$ws.Range("D21").Value = '4.807'
$ws.Range("E21").Value = '-7.09%'
$ws.Range("G21").Value = '3'
$ws.Range("E22").Value = '19.47%'
$ws.Range("G22").Value = '3'
$ws.Range("D23").Value = '0.04482'
$ws.Range("E23").Value = '1.59%'
$ws.Range("G23").Value = '3'
$ws.Range("D24").Value = '0.001244'
$ws.Range("E24").Value = '0.82%'
$ws.Range("G24").Value = '3'
$ws.Range("E25").Value = '-9.70%'
$ws.Range("G25").Value = '3'
$ws.Range("G26").Value = '3'
$ws.Range("D27").Value = '0.0001302'
$ws.Range("E27").Value = '0.00%'
$ws.Range("G27").Value = '3'
$ws.Range("G28").Value = '3'
$ws.Range("G29").Value = '3'
$ws.Range("G30").Value = '3'
$ws.Range("G31").Value = '3'
$ws.Range("G32").Value = '3'
$ws.Range("G33").Value = '3'
$ws.Range("G34").Value = '3'
$ws.Range("G35").Value = '3'
$ws.Range("G36").Value = '3'
$ws.Range("G37").Value = '3'
$ws.Range("G38").Value = '3'
$ws.Range("D39").Value = '0.02113'
$ws.Range("E39").Value = '9.05%'
$ws.Range("G39").Value = '3'
$ws.Range("D40").Value = '0.05149'
$ws.Range("E40").Value = '1.47%'
$ws.Range("G40").Value = '3'
$ws.Range("D41").Value = '0.007479'
$ws.Range("E41").Value = '-2.04%'
$ws.Range("G41").Value = '3'
$ws.Range("D42").Value = '0.01015'
$ws.Range("E42").Value = '-0.30%'
$ws.Range("G42").Value = '3'
$ws.Range("E43").Value = '1.76%'
$ws.Range("G43").Value = '3'
$ws.Range("D44").Value = '0.002133'
$ws.Range("E44").Value = '-1.39%'
$ws.Range("G44").Value = '3'
$ws.Range("D45").Value = '0.009685'
$ws.Range("E45").Value = '-0.92%'
$ws.Range("G45").Value = '3'
$ws.Range("D46").Value = '0.00006323'
$ws.Range("E46").Value = '2.63%'
$ws.Range("G46").Value = '3'
$ws.Range("E47").Value = '0.08%'
$ws.Range("G47").Value = '3'
$ws.Range("G48").Value = '3'
$ws.Range("E49").Value = '-3.47%'
$ws.Range("G49").Value = '3'
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").Value = '0.08%'
$ws.Range("G50").Value = '3'
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").Value = '0.08%'
$ws.Range("G51").Value = '3'

# Restore the default "Normal" style so we do not leave a stray text-format style behind
# on cells that originally had no explicit style.
$numericLookingRange.Style = "Normal"
